# counters_summary.xlsx - re-run of the data-quality scoring pass.
# The refreshed dataset zeroes out the previously-computed
# METADATACOMPLIANCE / COMPLETENESSOPTIONAL counters+scores for a few
# attributes, and appends a fresh (all-zero / blank-label) metadata row
# for the newly-introduced attribute at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Post unique reference: METADATACOMPLIANCE / METADATACOMPLIANCE SCORE reset
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Name: COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE reset
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# Job/Team Function: COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE reset
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# Contact phone: COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE reset
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

# Contact e-mail: COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE reset
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0

# Notes: COMPLETENESSOPTIONAL / COMPLETENESSOPTIONAL SCORE reset
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Append the generated metadata row for the new attribute: blank label,
# every counter/score column zeroed.
$ws.Range("A14").Value = ""
$ws.Range("B14:U14").Value = 0
